$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'nts''i''its'
$ws.Range("C2").Value = 'grandmother'
$ws.Range("E2").Value = 'noun'
$ws.Range("I2").Value = 'nts''iidz'
$ws.Range("B3").Value = 'niya''ay'
$ws.Range("C3").Value = 'grandfather'
$ws.Range("E3").Value = 'noun'
$ws.Range("B4").Value = 'nigwaat'
$ws.Range("C4").Value = 'father'
$ws.Range("E4").Value = 'noun'
$ws.Range("B5").Value = 'łgutx̱a''oo'
$ws.Range("C5").Value = 'cousin'
$ws.Range("E5").Value = 'noun'
$ws.Range("B6").Value = 'waky'
$ws.Range("C6").Value = 'a man''s brother'
$ws.Range("E6").Value = 'noun'
$ws.Range("B7").Value = 'insiipnsk'
$ws.Range("C7").Value = 'a woman''s brother'
$ws.Range("E7").Value = 'noun'
$ws.Range("B8").Value = 'hana''ax̱'
$ws.Range("C8").Value = 'woman'
$ws.Range("E8").Value = 'noun'
$ws.Range("B9").Value = 'łguułgu'
$ws.Range("C9").Value = 'child (of someone)'
$ws.Range("E9").Value = 'noun'
$ws.Range("B10").Value = 'kw''ida''ts'
$ws.Range("C10").Value = 'coat'
$ws.Range("E10").Value = 'noun'
$ws.Range("B11").Value = 'łimkt''ii'
$ws.Range("C11").Value = 'a man''s brother'
$ws.Range("E11").Value = 'noun'
$ws.Range("B12").Value = 'daala'
$ws.Range("C12").Value = 'dollar; money'
$ws.Range("E12").Value = 'noun'
$ws.Range("B13").Value = 'waa'
$ws.Range("C13").Value = 'name'
$ws.Range("E13").Value = 'noun'
$ws.Range("B14").Value = 'waalp'
$ws.Range("C14").Value = 'house'
$ws.Range("E14").Value = 'noun'
$ws.Range("I14").Value = 'waap'
$ws.Range("B15").Value = 'ts''ikts''ik'
$ws.Range("C15").Value = 'car'
$ws.Range("E15").Value = 'noun'
$ws.Range("B16").Value = 'x̱aldaawx̱k'
$ws.Range("C16").Value = 'medicine'
$ws.Range("E16").Value = 'noun'
$ws.Range("B17").Value = 'an''on'
$ws.Range("C17").Value = 'hand or arm'
$ws.Range("E17").Value = 'noun'
$ws.Range("B18").Value = 'naks'
$ws.Range("C18").Value = 'spouse'
$ws.Range("E18").Value = 'noun'
$ws.Range("B19").Value = 'kap'
$ws.Range("C19").Value = 'cup'
$ws.Range("E19").Value = 'noun'
$ws.Range("B20").Value = 'g̱oot'
$ws.Range("C20").Value = 'heart/mind'
$ws.Range("E20").Value = 'noun'
$ws.Range("B21").Value = 'wüliilm tgwah'
$ws.Range("C21").Value = 'glasses'
$ws.Range("E21").Value = 'noun'
$ws.Range("B22").Value = 'g̱aayt'
$ws.Range("C22").Value = 'hat'
$ws.Range("E22").Value = 'noun'
$ws.Range("B23").Value = 'hooya'
$ws.Range("C23").Value = 'clothes'
$ws.Range("E23").Value = 'noun'
$ws.Range("B24").Value = 'g̱aws'
$ws.Range("C24").Value = 'hair'
$ws.Range("E24").Value = 'noun'
$ws.Range("B25").Value = 'ha''ligyet'
$ws.Range("C25").Value = 'birthday'
$ws.Range("E25").Value = 'noun'

# Remove the two now-unused trailing rows (86 and 87)
$ws.Range("A86:A87").EntireRow.Delete()

# Update the active selection to match the edited workbook
$ws.Range("C26").Select()
